$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" / "Valor Mora" data between row 16 and row 18:
#   Row 16 was (1806, 19200) -> becomes (1808, 48000)
#   Row 17 stays (1807, 48000)
#   Row 18 was (1808, 48000) -> becomes (1806, 19200)
$ws.Range("E16").Value = "1808"
$ws.Range("F16").Value = 48000

$ws.Range("E18").Value = "1806"
$ws.Range("F18").Value = 19200
